$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental flag (was blank) -> "false", written as literal text (not the
# boolean FALSE). Stage it in a scratch cell with a leading apostrophe so
# Excel treats it as text, then copy/paste-values into B7 so the target
# cell keeps its original formatting (no quote-prefix style picked up).
$xlPasteValues = -4163
$scratch = $ws.Range("Z1")
$scratch.Value = "'false"
$scratch.Copy()
$ws.Range("B7").PasteSpecial($xlPasteValues)
$scratch.Clear()

# Date property refreshed to the new publication timestamp
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Description (was blank) -> populated with the ValueSet description
$ws.Range("B17").Value = "Levels of accumulated recovery debt from training or stress"
